$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts existing rows 4-14 down to 5-15
$ws.Rows.Item(4).Insert()

# Copy the (now shifted-down) row 5 formatting/values as the template for new row 4,
# since most columns repeat the same constant values across all data rows.
$ws.Range("A5:R5").Copy()
$ws.Range("A4:R4").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# Set the specific values that differ for the new record
$ws.Cells.Item(4, 4).Value = 44676   # D4 - Fecha
$ws.Cells.Item(4, 10).Value = 120    # J4 - Volumen
$ws.Cells.Item(4, 11).Value = 4000   # K4 - Precio minimo
$ws.Cells.Item(4, 12).Value = 4500   # L4 - Precio maximo
$ws.Cells.Item(4, 13).Value = 4250   # M4 - Precio promedio ponderado
$ws.Cells.Item(4, 16).Value = 71     # P4 - Precio $/Kg
